# Generate Report for Handback
# Row 3 in each sheet corresponds to the b07ffab4-5540-460d-9686-9f583923cf1a file,
# which has now been handed back (in sync with en-US) instead of merely "Ready for handoff".

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $statusHandedBack
$wsOverview.Range("C3").Value = $statusHandedBack

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $statusHandedBack
$wsZhCn.Range("H3").Value = "2016-03-18 14:38:35"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $statusHandedBack
$wsDeDe.Range("H3").Value = "2016-03-18 14:38:40"
